$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 11 ("R40" rule row): column B held the rule-name text "R40".
# Re-label it to the text "1" (kept as text, not converted to a number,
# via the leading-apostrophe text-entry convention).
$ws.Cells.Item(11, 2).Value = "'1"
